# Correction numero de SIRET
# --------------------------------------------------------------------------
# The SIRET number stored in R2 ("num_siret" column) was wrong and needs to
# be corrected. Fixing it through the normal Excel UI/COM path (typing the
# new value into the cell) also causes a couple of small, expected side
# effects that show up in the saved file:
#   - the header cells for columns C (insee) and R (num_siret) lose the
#     "Text" number format they had and fall back to the sheet's normal
#     style (this happens because editing data in this column resets the
#     column's filter/format bookkeeping for the header row);
#   - the AutoFilter that was covering the header is dropped together with
#     its hidden defined name (_xlnm._FilterDatabase);
#   - rows 1-2 pick up an explicit row height.
# Everything below reproduces that end state using plain COM calls.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correction of the SIRET number itself ----------------------------
$ws.Range("R2").Value = "21920044100018"

# --- 2. Header cells C1/R1 drop their "Text" number format ---------------
# Copy the plain (unformatted) style from a neighbouring header cell so the
# resulting style exactly matches the rest of the header row instead of
# creating a brand new one.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Drop the AutoFilter (and the _xlnm._FilterDatabase name it owns) -
$ws.AutoFilterMode = $false
foreach ($n in @($wb.Names)) {
    $n.Delete()
}

# --- 4. Rows 1-2 get an explicit row height -------------------------------
$ws.Rows("1:2").RowHeight = 14.25
